$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# parameters sheet: Elasticity value change
# ---------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("A2").Value = 0.02

# ---------------------------------------------------------------
# tech sheet: a few value tweaks + a new 0.0000 number format
# ---------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("tech")
$wsTech.Range("G3").Value = 0
$wsTech.Range("D4").Value = 8
$wsTech.Range("G4").NumberFormat = "0.0000"
$wsTech.Range("G4").Value = [double]"5.9999999999999995E-4"
$wsTech.Columns.Item(7).ColumnWidth = 6.5

# ---------------------------------------------------------------
# day_weights sheet: new weights
# ---------------------------------------------------------------
$wsDayWeights = $wb.Worksheets.Item("day_weights")
$wsDayWeights.Range("B2").Value = 199
$wsDayWeights.Range("B3").Value = 106
$wsDayWeights.Range("B4").Value = 61

# ---------------------------------------------------------------
# cap_factors sheet: recomputed hourly capacity factors
# ---------------------------------------------------------------
$wsCapFactors = $wb.Worksheets.Item("cap_factors")

$wsCapFactors.Range("G2").Value = [double]"4.0201005025125598E-3"
$wsCapFactors.Range("H2").Value = [double]"4.1170854271356763E-2"
$wsCapFactors.Range("I2").Value = [double]"0.15498492462311561"
$wsCapFactors.Range("J2").Value = [double]"0.34120100502512568"
$wsCapFactors.Range("K2").Value = [double]"0.51137185929648243"
$wsCapFactors.Range("L2").Value = [double]"0.63440201005025132"
$wsCapFactors.Range("M2").Value = [double]"0.70589447236180902"
$wsCapFactors.Range("N2").Value = [double]"0.7173165829145729"
$wsCapFactors.Range("O2").Value = [double]"0.67734170854271358"
$wsCapFactors.Range("P2").Value = [double]"0.57809547738693468"
$wsCapFactors.Range("Q2").Value = [double]"0.43412562814070349"
$wsCapFactors.Range("R2").Value = [double]"0.25252763819095481"
$wsCapFactors.Range("S2").Value = [double]"9.1608040201005048E-2"
$wsCapFactors.Range("T2").Value = [double]"1.7597989949748739E-2"
$wsCapFactors.Range("U2").Value = [double]"3.2160804020100472E-4"

$wsCapFactors.Range("G3").Value = [double]"2.1509433962264152E-3"
$wsCapFactors.Range("H3").Value = [double]"2.0773584905660381E-2"
$wsCapFactors.Range("I3").Value = [double]"7.4056603773584917E-2"
$wsCapFactors.Range("J3").Value = [double]"0.20321698113207551"
$wsCapFactors.Range("K3").Value = [double]"0.34733018867924531"
$wsCapFactors.Range("L3").Value = [double]"0.46702830188679251"
$wsCapFactors.Range("M3").Value = [double]"0.54499999999999993"
$wsCapFactors.Range("N3").Value = [double]"0.55942452830188683"
$wsCapFactors.Range("O3").Value = [double]"0.51176415094339622"
$wsCapFactors.Range("P3").Value = [double]"0.41073584905660382"
$wsCapFactors.Range("Q3").Value = [double]"0.27456603773584909"
$wsCapFactors.Range("R3").Value = [double]"0.13354716981132081"
$wsCapFactors.Range("S3").Value = [double]"3.995283018867922E-2"
$wsCapFactors.Range("T3").Value = [double]"8.6698113207547139E-3"
$wsCapFactors.Range("U3").Value = [double]"1.4150943396226421E-4"

$wsCapFactors.Range("G4").Value = [double]"1.8032786885245899E-3"
$wsCapFactors.Range("H4").Value = [double]"1.6032786885245912E-2"
$wsCapFactors.Range("I4").Value = [double]"5.4229508196721322E-2"
$wsCapFactors.Range("J4").Value = [double]"0.1274918032786885"
$wsCapFactors.Range("K4").Value = [double]"0.20950819672131141"
$wsCapFactors.Range("L4").Value = [double]"0.26809836065573772"
$wsCapFactors.Range("M4").Value = [double]"0.28267213114754092"
$wsCapFactors.Range("N4").Value = [double]"0.27595081967213131"
$wsCapFactors.Range("O4").Value = [double]"0.25488524590163941"
$wsCapFactors.Range("P4").Value = [double]"0.21809836065573759"
$wsCapFactors.Range("Q4").Value = 0.151
$wsCapFactors.Range("R4").Value = [double]"8.1081967213114725E-2"
$wsCapFactors.Range("S4").Value = [double]"3.1426229508196707E-2"
$wsCapFactors.Range("T4").Value = [double]"7.7049180327868824E-3"
$wsCapFactors.Range("U4").Value = [double]"9.8360655737705021E-5"

# ---------------------------------------------------------------
# Selection / active-cell bookkeeping to mirror the authored view state
# ---------------------------------------------------------------
$wsParams.Activate()
$wsParams.Range("C8").Select()

$wsTech.Activate()
$wsTech.Range("D4").Select()

$wsCapFactors.Activate()
$wsCapFactors.Range("B3:Y3").Select()

$wsDayWeights.Activate()
$wsDayWeights.Range("B2:B4").Select()
